$d = $word.ActiveDocument

# Collapse a range to the very end of the document body (after the last
# paragraph, before the sectPr) so the new paragraphs are appended there.
$end = $d.Content
$end.Collapse(0)

$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$newContent = '<w:p ' + $xmlNs + '/>' +
  '<w:p ' + $xmlNs + '><w:r><w:t>Project Approach:</w:t></w:r></w:p>' +
  '<w:p ' + $xmlNs + '>' +
    '<w:r><w:tab/><w:t xml:space="preserve">For this project, I approached the problem by first reading the requirements and then creating a general layout on paper. </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">I started by remembering the logic I used to create a prefix expression calculator in Java, last semester. I know we used popping from a stack and saving expressions, however in Haskell it would have to be different since it is a functional language. We cannot use variables </w:t></w:r>' +
    '<w:r><w:t>in the same way as in Java or Python</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> I started to break down the project into its main parts: Parsing the expression, Evaluate expression, Receive user input, and a main function.</w:t></w:r>' +
  '</w:p>' +
  '<w:p ' + $xmlNs + '/>' +
  '<w:p ' + $xmlNs + '><w:r><w:t>Project Organization:</w:t></w:r></w:p>'

$null = $end.InsertXML($newContent)
